# Fruta / hortaliza, semanal
# Re-shuffles the per-row date/volume/price data (columns D, M, N, O, P, S)
# across the existing data rows (2-41) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44809
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 27000
$ws.Cells.Item(2, 15).Value = 28000
$ws.Cells.Item(2, 16).Value = 27500
$ws.Cells.Item(2, 19).Value = 1375
$ws.Cells.Item(4, 4).Value = 44880
$ws.Cells.Item(4, 13).Value = 100
$ws.Cells.Item(4, 14).Value = 28000
$ws.Cells.Item(4, 15).Value = 30000
$ws.Cells.Item(4, 16).Value = 29000
$ws.Cells.Item(4, 19).Value = 1450
$ws.Cells.Item(5, 4).Value = 44879
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(5, 14).Value = 28000
$ws.Cells.Item(5, 15).Value = 30000
$ws.Cells.Item(5, 16).Value = 29000
$ws.Cells.Item(5, 19).Value = 1450
$ws.Cells.Item(6, 4).Value = 44810
$ws.Cells.Item(6, 14).Value = 27000
$ws.Cells.Item(6, 15).Value = 28000
$ws.Cells.Item(6, 16).Value = 27500
$ws.Cells.Item(6, 19).Value = 1375
$ws.Cells.Item(7, 4).Value = 44365
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 20000
$ws.Cells.Item(7, 15).Value = 21000
$ws.Cells.Item(7, 16).Value = 20500
$ws.Cells.Item(7, 19).Value = 1025
$ws.Cells.Item(8, 4).Value = 44784
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 27000
$ws.Cells.Item(8, 15).Value = 28000
$ws.Cells.Item(8, 16).Value = 27500
$ws.Cells.Item(8, 19).Value = 1375
$ws.Cells.Item(9, 4).Value = 44874
$ws.Cells.Item(9, 13).Value = 240
$ws.Cells.Item(9, 14).Value = 29000
$ws.Cells.Item(9, 15).Value = 30000
$ws.Cells.Item(9, 16).Value = 29500
$ws.Cells.Item(9, 19).Value = 1475
$ws.Cells.Item(11, 4).Value = 44315
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 20000
$ws.Cells.Item(11, 15).Value = 21000
$ws.Cells.Item(11, 16).Value = 20500
$ws.Cells.Item(11, 19).Value = 1025
$ws.Cells.Item(12, 4).Value = 44335
$ws.Cells.Item(12, 14).Value = 19000
$ws.Cells.Item(12, 15).Value = 20000
$ws.Cells.Item(12, 16).Value = 19500
$ws.Cells.Item(12, 19).Value = 975
$ws.Cells.Item(13, 4).Value = 44445
$ws.Cells.Item(13, 13).Value = 160
$ws.Cells.Item(13, 14).Value = 20000
$ws.Cells.Item(13, 15).Value = 21000
$ws.Cells.Item(13, 16).Value = 20500
$ws.Cells.Item(13, 19).Value = 1025
$ws.Cells.Item(14, 4).Value = 44882
$ws.Cells.Item(14, 13).Value = 120
$ws.Cells.Item(14, 14).Value = 28000
$ws.Cells.Item(14, 15).Value = 30000
$ws.Cells.Item(14, 16).Value = 29000
$ws.Cells.Item(14, 19).Value = 1450
$ws.Cells.Item(15, 4).Value = 44301
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 19000
$ws.Cells.Item(15, 16).Value = 18500
$ws.Cells.Item(15, 19).Value = 925
$ws.Cells.Item(16, 4).Value = 44473
$ws.Cells.Item(16, 13).Value = 40
$ws.Cells.Item(16, 14).Value = 19500
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 19750
$ws.Cells.Item(16, 19).Value = 988
$ws.Cells.Item(17, 4).Value = 44778
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 23000
$ws.Cells.Item(17, 15).Value = 24000
$ws.Cells.Item(17, 16).Value = 23500
$ws.Cells.Item(17, 19).Value = 1175
$ws.Cells.Item(18, 4).Value = 44466
$ws.Cells.Item(18, 13).Value = 100
$ws.Cells.Item(18, 14).Value = 20000
$ws.Cells.Item(18, 15).Value = 21000
$ws.Cells.Item(18, 16).Value = 20500
$ws.Cells.Item(18, 19).Value = 1025
$ws.Cells.Item(19, 4).Value = 44427
$ws.Cells.Item(19, 13).Value = 200
$ws.Cells.Item(19, 15).Value = 21000
$ws.Cells.Item(19, 16).Value = 20500
$ws.Cells.Item(19, 19).Value = 1025
$ws.Cells.Item(20, 4).Value = 44410
$ws.Cells.Item(20, 13).Value = 200
$ws.Cells.Item(20, 14).Value = 20000
$ws.Cells.Item(20, 15).Value = 21000
$ws.Cells.Item(20, 16).Value = 20500
$ws.Cells.Item(20, 19).Value = 1025
$ws.Cells.Item(21, 4).Value = 44467
$ws.Cells.Item(21, 13).Value = 200
$ws.Cells.Item(22, 4).Value = 44448
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 20000
$ws.Cells.Item(22, 15).Value = 21000
$ws.Cells.Item(22, 16).Value = 20500
$ws.Cells.Item(22, 19).Value = 1025
$ws.Cells.Item(23, 4).Value = 44350
$ws.Cells.Item(23, 14).Value = 19000
$ws.Cells.Item(23, 15).Value = 20000
$ws.Cells.Item(23, 16).Value = 19500
$ws.Cells.Item(23, 19).Value = 975
$ws.Cells.Item(24, 4).Value = 44434
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 20000
$ws.Cells.Item(24, 15).Value = 21000
$ws.Cells.Item(24, 16).Value = 20500
$ws.Cells.Item(24, 19).Value = 1025
$ws.Cells.Item(25, 4).Value = 44336
$ws.Cells.Item(25, 13).Value = 100
$ws.Cells.Item(25, 14).Value = 19500
$ws.Cells.Item(25, 15).Value = 20000
$ws.Cells.Item(25, 16).Value = 19750
$ws.Cells.Item(25, 19).Value = 988
$ws.Cells.Item(26, 4).Value = 44474
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 19000
$ws.Cells.Item(26, 15).Value = 20000
$ws.Cells.Item(26, 16).Value = 19500
$ws.Cells.Item(26, 19).Value = 975
$ws.Cells.Item(27, 4).Value = 44441
$ws.Cells.Item(27, 13).Value = 160
$ws.Cells.Item(28, 4).Value = 44782
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 23500
$ws.Cells.Item(28, 15).Value = 24000
$ws.Cells.Item(28, 16).Value = 23750
$ws.Cells.Item(28, 19).Value = 1188
$ws.Cells.Item(29, 4).Value = 44776
$ws.Cells.Item(29, 13).Value = 160
$ws.Cells.Item(30, 4).Value = 44418
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(31, 4).Value = 44442
$ws.Cells.Item(31, 13).Value = 140
$ws.Cells.Item(31, 14).Value = 20000
$ws.Cells.Item(31, 15).Value = 21000
$ws.Cells.Item(31, 16).Value = 20500
$ws.Cells.Item(31, 19).Value = 1025
$ws.Cells.Item(32, 4).Value = 44462
$ws.Cells.Item(32, 14).Value = 19500
$ws.Cells.Item(32, 15).Value = 20000
$ws.Cells.Item(32, 16).Value = 19750
$ws.Cells.Item(32, 19).Value = 988
$ws.Cells.Item(33, 4).Value = 44326
$ws.Cells.Item(33, 13).Value = 160
$ws.Cells.Item(33, 14).Value = 19500
$ws.Cells.Item(33, 15).Value = 20000
$ws.Cells.Item(33, 16).Value = 19750
$ws.Cells.Item(33, 19).Value = 988
$ws.Cells.Item(34, 4).Value = 44435
$ws.Cells.Item(34, 13).Value = 260
$ws.Cells.Item(34, 14).Value = 20000
$ws.Cells.Item(34, 15).Value = 22000
$ws.Cells.Item(34, 16).Value = 21115
$ws.Cells.Item(34, 19).Value = 1056
$ws.Cells.Item(35, 4).Value = 44781
$ws.Cells.Item(35, 14).Value = 23000
$ws.Cells.Item(35, 15).Value = 24000
$ws.Cells.Item(35, 16).Value = 23500
$ws.Cells.Item(35, 19).Value = 1175
$ws.Cells.Item(36, 4).Value = 44407
$ws.Cells.Item(36, 13).Value = 160
$ws.Cells.Item(36, 14).Value = 20000
$ws.Cells.Item(36, 15).Value = 21000
$ws.Cells.Item(36, 16).Value = 20500
$ws.Cells.Item(36, 19).Value = 1025
$ws.Cells.Item(37, 4).Value = 44431
$ws.Cells.Item(37, 14).Value = 21000
$ws.Cells.Item(37, 15).Value = 22000
$ws.Cells.Item(37, 16).Value = 21500
$ws.Cells.Item(37, 19).Value = 1075
$ws.Cells.Item(38, 4).Value = 44420
$ws.Cells.Item(38, 13).Value = 160
$ws.Cells.Item(38, 14).Value = 20000
$ws.Cells.Item(38, 15).Value = 21000
$ws.Cells.Item(38, 16).Value = 20500
$ws.Cells.Item(38, 19).Value = 1025
$ws.Cells.Item(39, 4).Value = 44343
$ws.Cells.Item(39, 14).Value = 19500
$ws.Cells.Item(39, 15).Value = 20000
$ws.Cells.Item(39, 16).Value = 19750
$ws.Cells.Item(39, 19).Value = 988
$ws.Cells.Item(40, 4).Value = 44417
$ws.Cells.Item(40, 14).Value = 20000
$ws.Cells.Item(40, 15).Value = 21000
$ws.Cells.Item(40, 16).Value = 20500
$ws.Cells.Item(40, 19).Value = 1025
$ws.Cells.Item(41, 4).Value = 44364
$ws.Cells.Item(41, 13).Value = 140
$ws.Cells.Item(41, 14).Value = 20000
$ws.Cells.Item(41, 15).Value = 21000
$ws.Cells.Item(41, 16).Value = 20500
$ws.Cells.Item(41, 19).Value = 1025
